# Auto-generated PowerShell COM-interop script.
# Applies the "Updated cryptos list ... with GitHub Actions" diff: refreshed
# Price/Volume(1h) figures for every coin row, plus a rank swap between
# "ImmutableX" (row 37) and "Maker" (row 38).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source file stores every Price cell as plain text (e.g. "215.30",
# "1.540"), not a number - so trailing zeros / thousands-style dot grouping
# survive. Excel auto-coerces a numeric-looking string typed into a
# General-formatted cell into a real number (dropping e.g. "215.30" ->
# 215.3), so pre-format the numeric-looking Price cells as Text before
# writing their new values, to preserve the original text formatting.
$ws.Range("D5:D6").NumberFormat = "@"
$ws.Range("D8:D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21:D28").NumberFormat = "@"
$ws.Range("D30:D33").NumberFormat = "@"
$ws.Range("D35:D36").NumberFormat = "@"
$ws.Range("D38:D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43:D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '25.924.66'
$ws.Range("E2").Value = '  +0.30%  '

# Row 3
$ws.Range("D3").Value = '1.641.10'
$ws.Range("E3").Value = '  +0.17%  '

# Row 4
$ws.Range("E4").Value = '  -0.08%  '

# Row 5
$ws.Range("D5").Value = '215.30'
$ws.Range("E5").Value = '  -0.05%  '

# Row 6
$ws.Range("D6").Value = '0.5082'
$ws.Range("E6").Value = '  +1.19%  '

# Row 7
$ws.Range("E7").Value = '  +0.05%  '

# Row 8
$ws.Range("D8").Value = '0.2562'
$ws.Range("E8").Value = '  -0.16%  '

# Row 9
$ws.Range("D9").Value = '0.06391'
$ws.Range("E9").Value = '  +0.12%  '

# Row 10
$ws.Range("D10").Value = '19.51'
$ws.Range("E10").Value = '  -0.97%  '

# Row 11
$ws.Range("D11").Value = '0.07765'
$ws.Range("E11").Value = '  +0.52%  '

# Row 12
$ws.Range("D12").Value = '4.301'
$ws.Range("E12").Value = '  +0.82%  '

# Row 13
$ws.Range("D13").Value = '1.648.67'
$ws.Range("E13").Value = '  +0.43%  '

# Row 14
$ws.Range("D14").Value = '0.5449'
$ws.Range("E14").Value = '  +0.05%  '

# Row 15
$ws.Range("D15").Value = '0.0₅7841'
$ws.Range("E15").Value = '  -0.85%  '

# Row 16
$ws.Range("D16").Value = '64.64'
$ws.Range("E16").Value = '  +1.34%  '

# Row 17
$ws.Range("D17").Value = '25.972.81'
$ws.Range("E17").Value = '  +0.42%  '

# Row 18
$ws.Range("E18").Value = '  +0.02%  '

# Row 19
$ws.Range("D19").Value = '197.60'
$ws.Range("E19").Value = '  -1.65%  '

# Row 20
$ws.Range("E20").Value = '  +1.86%  '

# Row 21
$ws.Range("D21").Value = '9.951'
$ws.Range("E21").Value = '  +0.37%  '

# Row 22
$ws.Range("D22").Value = '6.038'
$ws.Range("E22").Value = '  +1.28%  '

# Row 23
$ws.Range("D23").Value = '1.007'
$ws.Range("E23").Value = '  +0.22%  '

# Row 24
$ws.Range("D24").Value = '1.876'
$ws.Range("E24").Value = '  -2.18%  '

# Row 25
$ws.Range("D25").Value = '140.90'
$ws.Range("E25").Value = '  -0.40%  '

# Row 26
$ws.Range("D26").Value = '0.1144'
$ws.Range("E26").Value = '  +0.78%  '

# Row 27
$ws.Range("D27").Value = '6.882'
$ws.Range("E27").Value = '  +2.73%  '

# Row 28
$ws.Range("D28").Value = '15.72'
$ws.Range("E28").Value = '  +0.25%  '

# Row 29
$ws.Range("E29").Value = '  -0.62%  '

# Row 30
$ws.Range("D30").Value = '0.05023'
$ws.Range("E30").Value = '  +0.71%  '

# Row 31
$ws.Range("D31").Value = '3.262'
$ws.Range("E31").Value = '  -0.36%  '

# Row 32
$ws.Range("D32").Value = '3.181'
$ws.Range("E32").Value = '  -0.38%  '

# Row 33
$ws.Range("D33").Value = '1.540'
$ws.Range("E33").Value = '  +0.07%  '

# Row 34
$ws.Range("E34").Value = '  -0.54%  '

# Row 35
$ws.Range("D35").Value = '0.8933'
$ws.Range("E35").Value = '  +0.37%  '

# Row 36
$ws.Range("D36").Value = '2.588'

# Row 37
$ws.Range("B37").Value = 'Maker'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D37").Value = '1.127.71'
$ws.Range("E37").Value = '  -3.45%  '

# Row 38
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").Value = '0.5497'
$ws.Range("E38").Value = '  -1.46%  '

# Row 39
$ws.Range("D39").Value = '0.01555'
$ws.Range("E39").Value = '  -0.40%  '

# Row 40
$ws.Range("E40").Value = '  +0.05%  '

# Row 41
$ws.Range("D41").Value = '2.551'
$ws.Range("E41").Value = '  -0.46%  '

# Row 42
$ws.Range("D42").Value = '0.0₈131'
$ws.Range("E42").Value = '  +12.46%  '

# Row 43
$ws.Range("D43").Value = '5.627'
$ws.Range("E43").Value = '  -1.02%  '

# Row 44
$ws.Range("D44").Value = '0.8174'
$ws.Range("E44").Value = '  +1.50%  '

# Row 45
$ws.Range("D45").Value = '99.91'
$ws.Range("E45").Value = '  +0.27%  '

# Row 46
$ws.Range("D46").Value = '1.778.19'
$ws.Range("E46").Value = '  +0.35%  '

# Row 47
$ws.Range("D47").Value = '0.4525'
$ws.Range("E47").Value = '  -0.07%  '

# Row 48
$ws.Range("E48").Value = '  -0.29%  '

# Row 49
$ws.Range("E49").Value = '  -0.06%  '

# Row 50
$ws.Range("D50").Value = '0.05076'
$ws.Range("E50").Value = '  +0.11%  '

# Row 51
$ws.Range("E51").Value = '  +0.24%  '
